# Apply "Generate Report for Handback" changes
#
# 1) Overview sheet: Status text "In Translation" -> "Handed back: in sync with en-US"
#    for both language rows (zh-cn row -> E2/F2, de-de row -> E3/F3),
#    plus widen columns E and F.
# 2) zh-cn / de-de detail sheets:
#    - Fill in "Latest Target File" (col I) with the source file name (hyperlinked,
#      same target URL as column A's link for that row).
#    - Fill in "Latest Handback File" (col J) with the per-language handback xlf name.
#    - Update "Latest Handback DateTime" (col K) with a real timestamp.
#    - Widen columns C, I and J.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# The "Status" column (C) on the detail sheets shares the same underlying
# string as the Overview sheet's status cells ("In Translation"), so it also
# flips to the handed-back message.
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "725ce266-6b25-4e18-a4f5-2e52030ff621.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c331eaa1089c83e35471206ff9a77accb105ece7/e2e/725ce266-6b25-4e18-a4f5-2e52030ff621.md", "", "", "725ce266-6b25-4e18-a4f5-2e52030ff621.md") | Out-Null
$zh.Range("J2").Value = "725ce266-6b25-4e18-a4f5-2e52030ff621.1ebeeba6d7582000190cb04381112a3f9d7c12dd.zh-cn.xlf"

$zh.Range("I3").Value = "edc45554-3494-4771-90ca-6e7452958358.md"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c331eaa1089c83e35471206ff9a77accb105ece7/e2e/edc45554-3494-4771-90ca-6e7452958358.md", "", "", "edc45554-3494-4771-90ca-6e7452958358.md") | Out-Null
$zh.Range("J3").Value = "edc45554-3494-4771-90ca-6e7452958358.8df0f57e4c9280938aa60ef30349faeb497a91ea.zh-cn.xlf"

$zh.Range("K2").Value = "2016-09-07 14:59:31"
$zh.Range("K3").Value = "2016-09-07 14:59:31"

$zh.Columns.Item(3).ColumnWidth = 29.17
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "725ce266-6b25-4e18-a4f5-2e52030ff621.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c331eaa1089c83e35471206ff9a77accb105ece7/e2e/725ce266-6b25-4e18-a4f5-2e52030ff621.md", "", "", "725ce266-6b25-4e18-a4f5-2e52030ff621.md") | Out-Null
$de.Range("J2").Value = "725ce266-6b25-4e18-a4f5-2e52030ff621.1ebeeba6d7582000190cb04381112a3f9d7c12dd.de-de.xlf"

$de.Range("I3").Value = "edc45554-3494-4771-90ca-6e7452958358.md"
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c331eaa1089c83e35471206ff9a77accb105ece7/e2e/edc45554-3494-4771-90ca-6e7452958358.md", "", "", "edc45554-3494-4771-90ca-6e7452958358.md") | Out-Null
$de.Range("J3").Value = "edc45554-3494-4771-90ca-6e7452958358.8df0f57e4c9280938aa60ef30349faeb497a91ea.de-de.xlf"

$de.Range("K2").Value = "2016-09-07 14:59:48"
$de.Range("K3").Value = "2016-09-07 14:59:48"

$de.Columns.Item(3).ColumnWidth = 29.17
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17

Write-Output "done"
